$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A new price record was added to the daily log. It belongs right before the
# current row 169, so existing rows 169:252 shift down one position to
# 170:253 (and the sheet's used range grows from R252 to R253).
$ws.Rows(169).Insert()

# Populate the newly inserted row 169 with the new record's values.
$ws.Cells.Item(169, 1).Value = 4
$ws.Cells.Item(169, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(169, 3).Value = 'Los Lagos'
$ws.Cells.Item(169, 4).Value = 44466
$ws.Cells.Item(169, 5).Value = 10
$ws.Cells.Item(169, 6).Value = 100112006
$ws.Cells.Item(169, 7).Value = 'Repollo'
$ws.Cells.Item(169, 8).Value = 'Crespo record'
$ws.Cells.Item(169, 9).Value = 'Segunda'
$ws.Cells.Item(169, 10).Value = 750
$ws.Cells.Item(169, 11).Value = 1000
$ws.Cells.Item(169, 12).Value = 1000
$ws.Cells.Item(169, 13).Value = 1000
$ws.Cells.Item(169, 14).Value = '$/unidad'
$ws.Cells.Item(169, 15).Value = 'Región del Maule'
$ws.Cells.Item(169, 16).Value = 1000
$ws.Cells.Item(169, 17).Value = 1
$ws.Cells.Item(169, 18).Value = 'Hortaliza'
